$d = $word.ActiveDocument

function Find-ParagraphIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) Remove the four "Confused about watchlist..." / route questions
#    paragraphs (they were resolved / removed from the notes).
# ---------------------------------------------------------------------
$startIdx = Find-ParagraphIndex("Confused about watchlist index and routes… and action")
$endIdx = Find-ParagraphIndex("Should the route be nested under user?")
$pStart = $d.Paragraphs.Item($startIdx)
$pEnd = $d.Paragraphs.Item($endIdx)
$delRange = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$delRange.Delete()

# ---------------------------------------------------------------------
# 2) Insert a block of new Todo-list items right before the existing
#    "Be able to add and remove watchlists" item, and mark that
#    existing item as done by prefixing it with "X ".
# ---------------------------------------------------------------------
$anchorIdx = Find-ParagraphIndex("Be able to add and remove watchlists")
$anchor = $d.Paragraphs.Item($anchorIdx)

# List of (level, text) pairs to insert in order, level 0 = top ilvl,
# level 1 = indented ilvl.
$newItems = @(
    @{ Level = 0; Text = "Set up API fetch to get current stock price (use the Stock Price: Quote API)" },
    @{ Level = 1; Text = "Add to company view" },
    @{ Level = 1; Text = "Add to watchlist view?" },
    @{ Level = 0; Text = "Put labels and heading on charts" },
    @{ Level = 0; Text = "Set up default image on watchlist cards" },
    @{ Level = 0; Text = "Finish watchlist detail" },
    @{ Level = 1; Text = "Add functions to add and remove tickers" },
    @{ Level = 1; Text = "Add function to hide detail" },
    @{ Level = 1; Text = "Add function to show current stock price (future)" }
)

foreach ($item in $newItems) {
    $anchorRange = $anchor.Range
    $anchorRange.InsertParagraphBefore()
    $newIdx = $anchorIdx
    $newPara = $d.Paragraphs.Item($newIdx)
    $newPara.Range.Text = $item.Text
    if ($item.Level -eq 1) {
        $newPara.Range.ListFormat.ListIndent()
    }
    $anchorIdx = $anchorIdx + 1
    $anchor = $d.Paragraphs.Item($anchorIdx)
}

# Mark the original item as completed with an "X " prefix.
$anchor.Range.InsertBefore("X ")
